# Fruta / hortaliza, semanal
#
# The weekly refresh reshuffles which data row holds which record: column
# values in D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado), O (Origen) and P (Precio $/Kg) move
# between rows 2..21 while every other column (A,B,C,E,F,G,H,I,N,Q,R)
# stays constant. Snapshot the mutable columns for every row first, then
# write them back out according to the new row order so we never read a
# cell after it has already been overwritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 21

# Snapshot current (pre-edit) values for the columns that move as a unit.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $row = @{
        D = $ws.Cells.Item($r, 4).Value2
        J = $ws.Cells.Item($r, 10).Value2
        K = $ws.Cells.Item($r, 11).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        O = $ws.Cells.Item($r, 15).Value2
        P = $ws.Cells.Item($r, 16).Value2
    }
    $snapshot[$r] = $row
}

# destinationRow -> sourceRow (which row's pre-edit record now lands here)
$mapping = @{
    2  = 13
    3  = 14
    4  = 19
    5  = 15
    6  = 16
    7  = 2
    8  = 12
    9  = 3
    10 = 17
    11 = 5
    12 = 10
    13 = 21
    14 = 18
    15 = 4
    16 = 20
    17 = 9
    18 = 7
    19 = 11
    20 = 6
    21 = 8
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $data = $snapshot[$srcRow]

    $ws.Cells.Item($destRow, 4).Value2 = $data.D
    $ws.Cells.Item($destRow, 10).Value2 = $data.J
    $ws.Cells.Item($destRow, 11).Value2 = $data.K
    $ws.Cells.Item($destRow, 12).Value2 = $data.L
    $ws.Cells.Item($destRow, 13).Value2 = $data.M
    $ws.Cells.Item($destRow, 15).Value2 = $data.O
    $ws.Cells.Item($destRow, 16).Value2 = $data.P
}

Write-Host "Row reshuffle applied."
